$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.236.44"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.691.71"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "216.59"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "23.04"
$ws.Range("E8").Value = "  +13.36%  "
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("D10").Value = "0.0629"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.930.76"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "1.693.35"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("E15").Value = "  +4.96%  "
$ws.Range("D16").Value = "67.55"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "27.243.96"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "237.45"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "8.15"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("E23").Value = "  +5.13%  "
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").Value = "148.38"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.50"
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").Value = "1.573.30"
$ws.Range("E33").Value = "  +6.41%  "
$ws.Range("D34").Value = "3.25"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "0.954"
$ws.Range("E36").Value = "  +5.84%  "
$ws.Range("D37").Value = "0.606"
$ws.Range("E37").Value = "  +3.93%  "
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  +4.88%  "
$ws.Range("D41").Value = "69.56"
$ws.Range("E41").Value = "  +3.07%  "
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "2.27"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("D45").Value = "1.838.23"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "0.788"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "91.18"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("E48").Value = "  +5.85%  "
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("E50").Value = "  +3.16%  "
$ws.Range("E51").Value = "  +6.29%  "
